$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 382, shifting existing rows 382-402 down to 383-403
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with the new data record
$ws.Cells.Item(382, 1).Value = 9
$ws.Cells.Item(382, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(382, 3).Value = "Metropolitana"
$ws.Cells.Item(382, 4).Value = 44783
$ws.Cells.Item(382, 5).Value = 13
$ws.Cells.Item(382, 6).Value = 100112052
$ws.Cells.Item(382, 7).Value = "Albahaca"
$ws.Cells.Item(382, 8).Value = "Sin especificar"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 200
$ws.Cells.Item(382, 11).Value = 5000
$ws.Cells.Item(382, 12).Value = 5000
$ws.Cells.Item(382, 13).Value = 5000
$ws.Cells.Item(382, 14).Value = "$/paquete"
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 5000
$ws.Cells.Item(382, 17).Value = 1
$ws.Cells.Item(382, 18).Value = "Hortaliza"
